$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a D-column (Price) cell to stay plain text even when the
# string looks like a number (e.g. "311.71"), mirroring the source data
# where every Price/Volume cell is stored as text.
function Set-TextValue($range, $value) {
    if ($value -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $range.NumberFormat = "@"
    }
    $range.Value = $value
}

# --- Row 30/31 content swap: InjectiveProtocol <-> Monero (ranking reorder) ---
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D30") '166.08'
$ws.Range("E30").Value = '  -4.57%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D31") '35.69'
$ws.Range("E31").Value = '  +0.56%  '

# --- Row 49/50 content swap: ordi <-> FraxShare (ranking reorder) ---
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D49") '9.00'
$ws.Range("E49").Value = '  -0.21%  '
$ws.Range("B50").Value = 'ordi'
$ws.Range("C50").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextValue $ws.Range("D50") '75.26'
$ws.Range("E50").Value = '  +10.58%  '

# --- Per-row Price (D) / Volume(1h) (E) refresh ---
Set-TextValue $ws.Range("D2") '42.409.24'
$ws.Range("E2").Value = '  -0.12%  '
Set-TextValue $ws.Range("D3") '2.328.02'
$ws.Range("E3").Value = '  -1.36%  '
$ws.Range("E4").Value = '  +0.12%  '
Set-TextValue $ws.Range("D5") '311.71'
$ws.Range("E5").Value = '  -4.78%  '
Set-TextValue $ws.Range("D6") '106.17'
$ws.Range("E6").Value = '  +5.89%  '
Set-TextValue $ws.Range("D7") '0.629'
$ws.Range("E7").Value = '  -1.17%  '
$ws.Range("E8").Value = '  +0.14%  '
Set-TextValue $ws.Range("D9") '0.606'
$ws.Range("E9").Value = '  -2.30%  '
Set-TextValue $ws.Range("D10") '39.95'
$ws.Range("E10").Value = '  -0.49%  '
Set-TextValue $ws.Range("D11") '0.0918'
$ws.Range("E11").Value = '  -0.06%  '
Set-TextValue $ws.Range("D12") '8.38'
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("E13").Value = '  +0.66%  '
Set-TextValue $ws.Range("D14") '0.981'
$ws.Range("E14").Value = '  -3.07%  '
Set-TextValue $ws.Range("D15") '15.65'
$ws.Range("E15").Value = '  -3.51%  '
Set-TextValue $ws.Range("D16") '2.680.70'
$ws.Range("E16").Value = '  -1.50%  '
Set-TextValue $ws.Range("D17") '2.352.64'
$ws.Range("E17").Value = '  -0.49%  '
Set-TextValue $ws.Range("D18") '42.384.00'
$ws.Range("E18").Value = '  -0.42%  '
$ws.Range("E19").Value = '  -0.69%  '
$ws.Range("E20").Value = '  -0.97%  '
Set-TextValue $ws.Range("D21") '75.61'
$ws.Range("E21").Value = '  +0.66%  '
Set-TextValue $ws.Range("D22") '3.50'
$ws.Range("E22").Value = '  -5.45%  '
Set-TextValue $ws.Range("D23") '265.13'
$ws.Range("E23").Value = '  -3.74%  '
Set-TextValue $ws.Range("D24") '2.29'
$ws.Range("E24").Value = '  -0.43%  '
Set-TextValue $ws.Range("D25") '9.32'
$ws.Range("E25").Value = '  -3.81%  '
$ws.Range("E26").Value = '  +0.42%  '
Set-TextValue $ws.Range("D27") '11.13'
$ws.Range("E27").Value = '  -3.02%  '
Set-TextValue $ws.Range("D28") '23.18'
$ws.Range("E28").Value = '  -2.70%  '
$ws.Range("E29").Value = '  +1.71%  '
$ws.Range("E32").Value = '  -0.46%  '
Set-TextValue $ws.Range("D33") '2.91'
$ws.Range("E33").Value = '  -6.32%  '
Set-TextValue $ws.Range("D34") '5.95'
$ws.Range("E34").Value = '  +0.13%  '
$ws.Range("E35").Value = '  +14.23%  '
$ws.Range("E36").Value = '  -2.94%  '
Set-TextValue $ws.Range("D37") '4.56'
$ws.Range("E37").Value = '  -0.97%  '
Set-TextValue $ws.Range("D38") '0.0354'
$ws.Range("E38").Value = '  -1.09%  '
Set-TextValue $ws.Range("D39") '3.75'
$ws.Range("E39").Value = '  -3.37%  '
Set-TextValue $ws.Range("D40") '2.62'
$ws.Range("E40").Value = '  -8.80%  '
Set-TextValue $ws.Range("D41") '104.71'
$ws.Range("E41").Value = '  +15.81%  '
Set-TextValue $ws.Range("D42") '0.234'
$ws.Range("E42").Value = '  +2.79%  '
$ws.Range("E43").Value = '  -3.27%  '
Set-TextValue $ws.Range("D44") '70.57'
$ws.Range("E44").Value = '  +2.53%  '
$ws.Range("E45").Value = '  +0.10%  '
Set-TextValue $ws.Range("D46") '12.14'
$ws.Range("E46").Value = '  +1.98%  '
Set-TextValue $ws.Range("D47") '111.88'
$ws.Range("E47").Value = '  -3.14%  '
Set-TextValue $ws.Range("D48") '5.43'
$ws.Range("E48").Value = '  -0.32%  '
Set-TextValue $ws.Range("D51") '1.26'
$ws.Range("E51").Value = '  -0.06%  '
